# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values for rows 2-47, replacing the old Strike# based values.
$kValues = @(1,0,2,1,1,3,2,1,3,0,1,2,1,0,0,4,0,3,2,3,1,0,4,1,3,0,1,0,3,1,3,1,0,0,1,3,3,0,3,3,3,1,3,0,3,2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
